$wb = $excel.ActiveWorkbook

function Update-Value {
    param($ws, $cellRef, $newValue)
    $ws.Range($cellRef).Value = $newValue
}

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
Update-Value $ws1 "F3" 1195
Update-Value $ws1 "F4" 16811
Update-Value $ws1 "F6" 1645
Update-Value $ws1 "F7" 67
Update-Value $ws1 "F9" 382
Update-Value $ws1 "F10" 218
Update-Value $ws1 "F12" 11655
Update-Value $ws1 "F14" 1335
Update-Value $ws1 "F15" 4623
Update-Value $ws1 "F16" 451
Update-Value $ws1 "F18" 68
Update-Value $ws1 "F19" 894

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
Update-Value $ws4 "F4" 1195
Update-Value $ws4 "F5" 16811
Update-Value $ws4 "F7" 1645
Update-Value $ws4 "F8" 67
Update-Value $ws4 "F10" 382
Update-Value $ws4 "F11" 218
Update-Value $ws4 "F15" 11655
Update-Value $ws4 "F17" 1335
Update-Value $ws4 "F18" 4623
Update-Value $ws4 "F19" 451
Update-Value $ws4 "F21" 68
Update-Value $ws4 "F22" 894
